# ---------------------------------------------------------------------------
# Fills in the two blank "matching" exercise templates on sheets "8_" and
# "9_" with new term/definition content, widens/adds the 4th (comment)
# column, extends each table by the needed rows, and finally leaves the
# workbook's active sheet/selection state the way the author left it
# (sheet "9_" active, plus the exact cell selections recorded on sheets
# "8_", "9_" and "Matching").
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ===========================================================================
# Sheet "8_"  (internal worksheets/sheet9.xml)
# ===========================================================================
$ws8 = $wb.Worksheets.Item("8_")

# --- text content (columns A-C already exist with the right fill/wrap style
#     from the template; just overwrite the values) ------------------------
$ws8.Range("A1").Value = "Match each situation ('term') with the word that best describes that situation ('definition')."
$ws8.Range("B1").Value = "Correct order of definitions"
$ws8.Range("C1").Value = "Definitions"

$ws8.Range("A2").Value = "The government lists the detailed steps that must be taken every time a brake pad is tested"
$ws8.Range("B2").Value = "B"
$ws8.Range("C2").Value = "Code"

$ws8.Range("A3").Value = "The government sets a maximum temperature that a brake pad can reach during an independent test in order to make that pad available to the public"
$ws8.Range("B3").Value = "D"
$ws8.Range("C3").Value = "Protocol"

$ws8.Range("A4").Value = "A brake manufacturer sets a lower maximum temperature that its own brake pad can reach during its corporate tests"
$ws8.Range("B4").Value = "C"
$ws8.Range("C4").Value = "Internal Standard"

$ws8.Range("A5").Value = "The government releases a set of rules defining the minimum requirements for brake pads "
$ws8.Range("B5").Value = "A"
$ws8.Range("C5").Value = "Regulatory Standard"

# --- new column D: clone formats from column C (wrap, same fill) then fill
# in the ones that actually carry a comment ---------------------------------
$ws8.Range("C1:C5").Copy()
$ws8.Range("D1:D5").PasteSpecial($xlPasteFormats)

$ws8.Range("D2").Value = "A protocol usually describes a set of steps."
$ws8.Range("D4").Value = "An internal standard is generally more stringent than a regulatory standard."
$ws8.Range("D5").Value = "A code usually refers to a set of rules, rather than a single rule."

# --- new padding rows 6 & 7 -------------------------------------------------
# Columns A & D keep the normal wrapped data style (same as column A rows 2-5)
$ws8.Range("A5").Copy()
$ws8.Range("A6:A7").PasteSpecial($xlPasteFormats)
$ws8.Range("D5").Copy()
$ws8.Range("D6:D7").PasteSpecial($xlPasteFormats)

# Columns B & C on the padding rows use the unwrapped "Matching" filler style
$wsMatching = $wb.Worksheets.Item("Matching")
$wsMatching.Range("A6").Copy()
$ws8.Range("B6:C7").PasteSpecial($xlPasteFormats)

# --- column widths & row heights -------------------------------------------
$ws8.Columns.Item(3).ColumnWidth = 21.833333333333336   # -> stored width 22.7109375
$ws8.Columns.Item(4).ColumnWidth = 34.166666666666664   # -> stored width 35

$ws8.Rows.Item(1).RowHeight = 45
$ws8.Rows.Item(2).RowHeight = 45
$ws8.Rows.Item(3).RowHeight = 75
$ws8.Rows.Item(4).RowHeight = 60
$ws8.Rows.Item(5).RowHeight = 45

# ===========================================================================
# Sheet "9_"  (internal worksheets/sheet10.xml)
# ===========================================================================
$ws9 = $wb.Worksheets.Item("9_")

$ws9.Range("A1").Value = "Match the desirable behavior/ outcome with justification for using standards"
$ws9.Range("B1").Value = "Correct"
$ws9.Range("C1").Value = "Comment"

$ws9.Range("A2").Value = "A customer reads that a package of cookies has 180 calories, and can dependably know that this package has fewer calories than the package that says 220 calories"
$ws9.Range("B2").Value = "D"
$ws9.Range("C2").Value = "Liability "

$ws9.Range("A3").Value = "A exit door explodes from an airplane, but you know when your team decided which bolts to use, they were following regulatory standards"
$ws9.Range("B3").Value = "A"
$ws9.Range("C3").Value = "Safety"

$ws9.Range("A4").Value = "A city hired a contractor to build a bridge, but didn't open the bridge when an inspector noticed that the beams used were not up to the standard"
$ws9.Range("B4").Value = "B"
$ws9.Range("C4").Value = "Interoperability"

$ws9.Range("A5").Value = "You bought an off-brand charger for your iPhone, but it worked without any problems"
$ws9.Range("B5").Value = "C"
$ws9.Range("C5").Value = "Communication"

# --- new column D -----------------------------------------------------------
$ws9.Range("C1:C5").Copy()
$ws9.Range("D1:D5").PasteSpecial($xlPasteFormats)

$ws9.Range("D2").Value = "Using the same processes to test and describe products makes it easier to compare those products."
$ws9.Range("D3").Value = "Standards can provide protection against lawsuits in the event of product failure."
$ws9.Range("D4").Value = "Enforcing standards can keep unsafe products out of the marketplace."
$ws9.Range("D5").Value = "Interoperability can be helpful to customers and manufacturers."

# --- new row 6 (a full matching row) & blank padding row 7 ------------------
$ws9.Range("A5:D5").Copy()
$ws9.Range("A6:D6").PasteSpecial($xlPasteFormats)
$ws9.Range("A7:D7").PasteSpecial($xlPasteFormats)

$ws9.Range("A6").Value = "A customer bought your company's refrigerator, which met the stringent internal standard for noise production, and was pleased with how quiet it was"
$ws9.Range("B6").Value = "E"
$ws9.Range("C6").Value = "Reputation/ public trust"
$ws9.Range("D6").Value = "A brand can gain customer loyalty (and charge more) if customers have good experiences with their products."

# --- column widths & row heights --------------------------------------------
$ws9.Columns.Item(1).ColumnWidth = 44.166666666666664   # -> stored width 45
$ws9.Columns.Item(3).ColumnWidth = 17.333333333333336   # -> stored width 18.140625
$ws9.Columns.Item(4).ColumnWidth = 41.33333333333333    # -> stored width 42.140625

$ws9.Rows.Item(1).RowHeight = 30
$ws9.Rows.Item(2).RowHeight = 60
$ws9.Rows.Item(3).RowHeight = 45
$ws9.Rows.Item(4).RowHeight = 60
$ws9.Rows.Item(5).RowHeight = 30
$ws9.Rows.Item(6).RowHeight = 60

# ===========================================================================
# Selections / active sheet, matching the author's final view state
# ===========================================================================
$ws8.Activate()
$ws8.Range("A2:D7").Select()

$wsMatching.Activate()
$wsMatching.Range("B1:D7").Select()

$ws9.Activate()
$ws9.Range("C14").Select()
